# Insert a new row at position 125 (shifting existing rows 125..220 down to 126..221)
# and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(125).Insert()

$ws.Cells.Item(125, 1).Value = 8
$ws.Cells.Item(125, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(125, 3).Value = "Coquimbo"
$ws.Cells.Item(125, 4).Value = "2023-11-06"
$ws.Cells.Item(125, 5).Value = 4
$ws.Cells.Item(125, 6).Value = 100112052
$ws.Cells.Item(125, 7).Value = "Albahaca"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 1000
$ws.Cells.Item(125, 11).Value = 3300
$ws.Cells.Item(125, 12).Value = 3500
$ws.Cells.Item(125, 13).Value = 3400
$ws.Cells.Item(125, 14).Value = "`$/paquete"
$ws.Cells.Item(125, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(125, 16).Value = 3400
$ws.Cells.Item(125, 17).Value = 1
$ws.Cells.Item(125, 18).Value = "Hortaliza"
